$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.673.13'
$ws.Range('E2').Value = '  +3.97%  '
$ws.Range('D3').Value = '3.021.85'
$ws.Range('E3').Value = '  +3.45%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '564.76'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.24%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.87'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +8.26%  '
$ws.Range('E8').Value = '  +2.25%  '
$ws.Range('D9').Value = '3.012.94'
$ws.Range('E9').Value = '  +3.51%  '
$ws.Range('E10').Value = '  +6.49%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.27'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +11.19%  '
$ws.Range('E12').Value = '  +3.37%  '
$ws.Range('E13').Value = '  +6.14%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.99%  '
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').Value = '3.526.33'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.24'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +6.40%  '
$ws.Range('D18').Value = '3.020.50'
$ws.Range('E18').Value = '  +3.86%  '
$ws.Range('D19').Value = '59.681.52'
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '436.48'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.90%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.68'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.44%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.724'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +6.53%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.14'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '13.29'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.22%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '80.85'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.26'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +14.44%  '
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +3.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.84'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +5.80%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.06'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.62%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.29'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +6.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.102'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.58%  '
$ws.Range('D34').Value = '0.0₃0785'
$ws.Range('E34').Value = '  +16.57%  '
$ws.Range('E35').Value = '  +7.92%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.95'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.85%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.12'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('E38').Value = '  +3.10%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.68'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.80'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +10.39%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '403.08'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.91%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0355'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.51%  '
$ws.Range('D43').Value = '2.786.73'
$ws.Range('E43').Value = '  +4.88%  '
$ws.Range('E44').Value = '  +0.71%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.255'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +7.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '123.26'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('E48').Value = '  +1.96%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.02'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.46%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.70'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +21.54%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '23.56'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.17%  '
